# edit.ps1 — "Writing a new Article for Perfecting the Fight Room"
#
# Four content changes to the SEO meta "header" document:
#   1. <meta name="title">            13 Door System   -> 18 Perfecting the Fight
#   2. <meta name="description">      rewritten, with the new sentence wrapped
#                                      in a Word "_Hlk222220365" bookmark
#   3. <meta name="revised">          January, 12        -> February 17
#   4. <meta name="url">              13_Door_System/... -> 18_Perfecting_The_Fight/19_Perfecting_the_Fight.html
#
# Note on technique: Find/Replace (like Word itself) re-coalesces every run
# in the touched paragraph that shares identical formatting, which can
# merge runs that the target markup keeps distinct (they only differ by
# Word's internal rsid bookkeeping, not by anything visible). We restore
# those run boundaries the same way Word does when you simply click/select
# at that spot: drop a zero-length bookmark there (forces a run split) and
# immediately remove the bookmark again — the split stays behind, nothing
# else does.

$d = $word.ActiveDocument

$script:splitCounter = 0
function Split-RunAt($doc, $pos) {
    $script:splitCounter = $script:splitCounter + 1
    $tmpName = "TmpSplitMarker" + $script:splitCounter
    $r = $doc.Range($pos, $pos)
    $doc.Bookmarks.Add($tmpName, $r) | Out-Null
    $doc.Bookmarks($tmpName).Delete()
}

# ---------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------
$oldTitle = "13 Door System"
$newTitle = "18 Perfecting the Fight"
$d.Content.Find.Execute($oldTitle, $false, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2) | Out-Null

$titlePara = $d.Paragraphs(1)
$pStart = $titlePara.Range.Start
$relStart = $titlePara.Range.Text.IndexOf($newTitle)
$absStart = $pStart + $relStart
$absEnd = $absStart + $newTitle.Length
Split-RunAt $d $absStart
Split-RunAt $d $absEnd

# ---------------------------------------------------------------------
# 2. Description — replace the sentence and wrap it in a bookmark, same
#    as Word leaves behind when that smart-paste/citation tracking fires.
# ---------------------------------------------------------------------
$oldDesc = "In this tutorial, we will be starting to build our door system. In the Dungeon Crawler game, we will want to have a variety of different doors, which can be coded to take the hero to various areas, in the game."
$descPart1 = "In this tutorial, "
$descPart2 = "we will be going back into our fight room object code and making a few changes, to create a more convincing fighting illusion, between the two characters."
$newDesc = $descPart1 + $descPart2 + " "
$d.Content.Find.Execute($oldDesc, $false, $false, $false, $false, $false, $true, 1, $false, $newDesc, 2) | Out-Null

$descPara = $d.Paragraphs(9)
$pStart = $descPara.Range.Start
$relStart = $descPara.Range.Text.IndexOf($descPart1)
$bmStart = $pStart + $relStart
$bmEnd = $bmStart + ($descPart1 + $descPart2).Length

# Wrap the new sentence in the bookmark Word left behind.
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_Hlk222220365", $bmRange) | Out-Null

# Split "In this tutorial, " from the rest (inside the bookmark).
$splitPoint1 = $bmStart + $descPart1.Length
Split-RunAt $d $splitPoint1

# Split the trailing " " from the closing `"/>` run (just after bookmark end).
$splitPoint2 = $bmEnd + 1
Split-RunAt $d $splitPoint2

# ---------------------------------------------------------------------
# 3. Revised date
# ---------------------------------------------------------------------
$oldDate = "January, 12"
$newDate = "February 17"
$d.Content.Find.Execute($oldDate, $false, $false, $false, $false, $false, $true, 1, $false, $newDate, 2) | Out-Null

$datePara = $d.Paragraphs(18)
$pStart = $datePara.Range.Start
$relStart = $datePara.Range.Text.IndexOf($newDate)
$absStart = $pStart + $relStart
$absEnd = $absStart + $newDate.Length
Split-RunAt $d $absStart
Split-RunAt $d $absEnd

# The trailing ", 2026" run gets coalesced with the closing `" />` run by
# the same re-coalescing; split those back apart too.
$relStart2 = $datePara.Range.Text.IndexOf('" /')
$splitPoint3 = $pStart + $relStart2
Split-RunAt $d $splitPoint3

# ---------------------------------------------------------------------
# 4. URL
# ---------------------------------------------------------------------
$oldUrl = "Enlightenment/Articles/2026/2_Game_Maker_2/13_Door_System/13_Door_System.html"
$newUrl = "Enlightenment/Articles/2026/2_Game_Maker_2/18_Perfecting_The_Fight/19_Perfecting_the_Fight.html"
$d.Content.Find.Execute($oldUrl, $false, $false, $false, $false, $false, $true, 1, $false, $newUrl, 2) | Out-Null

$urlPara = $d.Paragraphs(20)
$pStart = $urlPara.Range.Start
$relStart = $urlPara.Range.Text.IndexOf($newUrl)
$absStart = $pStart + $relStart
$absEnd = $absStart + $newUrl.Length
Split-RunAt $d $absStart
Split-RunAt $d $absEnd

Write-Host "Done."
